$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column D (shifts existing quarterly data from D:K to F:M)
$ws.Columns("D:E").Insert()

# 2. Copy number formatting from the (just-shifted) F:G columns into the new D:E columns
#    so the new quarter columns inherit the same date/number formats as the rest of the table.
$ws.Range("F7:G102").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 3. Match column widths: new D:E (newest two quarters) should use the "date column" width,
#    same width that column H (shifted from the old F) already has.
$ws.Columns("D:E").ColumnWidth = $ws.Columns("H").ColumnWidth

# 4. Populate the two newly-inserted columns (D, E) with the new quarter figures
$ws.Range("D7").Value2 = 43464
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 756300
$ws.Range("E8").Value2 = 674300
$ws.Range("D9").Value2 = 380100
$ws.Range("E9").Value2 = 342000
$ws.Range("D10").Value2 = 376200
$ws.Range("E10").Value2 = 332300
$ws.Range("D12").Value2 = 52000
$ws.Range("E12").Value2 = 48800
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = -1800
$ws.Range("E14").Value2 = -6500
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 640900
$ws.Range("E17").Value2 = 581100
$ws.Range("D18").Value2 = 115400
$ws.Range("E18").Value2 = 93200
$ws.Range("D20").Value2 = -19800
$ws.Range("E20").Value2 = 1500
$ws.Range("D21").Value2 = 142900
$ws.Range("E21").Value2 = 139900
$ws.Range("D22").Value2 = 16200
$ws.Range("E22").Value2 = 16700
$ws.Range("D23").Value2 = 79400
$ws.Range("E23").Value2 = 78000
$ws.Range("D24").Value2 = 10100
$ws.Range("E24").Value2 = 2600
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 69300
$ws.Range("E26").Value2 = 75400
$ws.Range("D27").Value2 = 69300
$ws.Range("E27").Value2 = 75400
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 2000
$ws.Range("E29").Value2 = 1100
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 19800
$ws.Range("E32").Value2 = -1500
$ws.Range("D33").Value2 = 71300
$ws.Range("E33").Value2 = 76500
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 71300
$ws.Range("E35").Value2 = 76500
$ws.Range("D38").Value2 = 43464
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 163100
$ws.Range("E41").Value2 = 149500
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 632700
$ws.Range("E43").Value2 = 551400
$ws.Range("D44").Value2 = 338300
$ws.Range("E44").Value2 = 354200
$ws.Range("D45").Value2 = 100500
$ws.Range("E45").Value2 = 110400
$ws.Range("D46").Value2 = 1234600
$ws.Range("E46").Value2 = 1165500
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 318600
$ws.Range("E48").Value2 = 317000
$ws.Range("D49").Value2 = 4152300
$ws.Range("E49").Value2 = 4144500
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 270000
$ws.Range("E52").Value2 = 235300
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 5975500
$ws.Range("E54").Value2 = 5862400
$ws.Range("D57").Value2 = 220900
$ws.Range("E57").Value2 = 180700
$ws.Range("D58").Value2 = 14900
$ws.Range("E58").Value2 = 20100
$ws.Range("D59").Value2 = 535800
$ws.Range("E59").Value2 = 486000
$ws.Range("D60").Value2 = 771600
$ws.Range("E60").Value2 = 686800
$ws.Range("D61").Value2 = 1876600
$ws.Range("E61").Value2 = 1882500
$ws.Range("D62").Value2 = 742300
$ws.Range("E62").Value2 = 720600
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 3390600
$ws.Range("E66").Value2 = 3289900
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 2602100
$ws.Range("E72").Value2 = 2532100
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 2585000
$ws.Range("E76").Value2 = 2572500
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43464
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 71300
$ws.Range("E81").Value2 = 76500
$ws.Range("D83").Value2 = 47200
$ws.Range("E83").Value2 = 45200
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 159400
$ws.Range("E89").Value2 = 93200
$ws.Range("D91").Value2 = -32800
$ws.Range("E91").Value2 = -20800
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -88000
$ws.Range("E94").Value2 = 8000
$ws.Range("D96").Value2 = -7800
$ws.Range("E96").Value2 = -7800
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -57300
$ws.Range("E100").Value2 = -109000
$ws.Range("D101").Value2 = -600
$ws.Range("E101").Value2 = -3100
$ws.Range("D102").Value2 = 13600
$ws.Range("E102").Value2 = -10900

# 5. A handful of older quarters were restated with slightly different figures;
#    apply those corrections on top of the shifted values.
$ws.Range("H9").Value2 = 334200
$ws.Range("I9").Value2 = 285300
$ws.Range("H10").Value2 = 307400
$ws.Range("I10").Value2 = 269000
$ws.Range("H17").Value2 = 548100
$ws.Range("I17").Value2 = 476200
$ws.Range("H18").Value2 = 93500
$ws.Range("I18").Value2 = 78100
$ws.Range("H20").Value2 = -1200
$ws.Range("I20").Value2 = 38000
$ws.Range("H32").Value2 = 1200
$ws.Range("I32").Value2 = -38000
$ws.Range("F47").Value2 = 0
$ws.Range("G47").Value2 = 0
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("J47").Value2 = 0
$ws.Range("I91").Value2 = -10900
$ws.Range("J91").Value2 = -5500
$ws.Range("I101").Value2 = 5000

